$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Pre-format the new K:P block as Text + left/center aligned.
#    This mirrors the author applying a text number-format + left alignment
#    to the whole new block before typing values, creating a new cellXfs entry
#    (numFmtId 49 "@", applyNumberFormat + applyAlignment) and widening the
#    sheet dimension / row spans out to column P.
$block = $ws.Range("K1:P42")
$block.NumberFormat = "@"
$block.HorizontalAlignment = -4131
$block.VerticalAlignment = -4108

# 2) Fill in every populated cell of the new MIPS op/rs/rt/rd/shamt/func table.
$ws.Range("K1").Value = "op"
$ws.Range("L1").Value = "rs"
$ws.Range("M1").Value = "rt"
$ws.Range("N1").Value = "rd"
$ws.Range("O1").Value = "shamt"
$ws.Range("P1").Value = "func"
$ws.Range("K2").Value = "000000"
$ws.Range("O2").Value = "00000"
$ws.Range("P2").Value = "100000"
$ws.Range("K3").Value = "000000"
$ws.Range("O3").Value = "00000"
$ws.Range("P3").Value = "100001"
$ws.Range("K4").Value = "001000"
$ws.Range("K5").Value = "001001"
$ws.Range("K6").Value = "000000"
$ws.Range("O6").Value = "00000"
$ws.Range("P6").Value = "100010"
$ws.Range("K7").Value = "000000"
$ws.Range("O7").Value = "00000"
$ws.Range("P7").Value = "100011"
$ws.Range("K8").Value = "011100"
$ws.Range("O8").Value = "00000"
$ws.Range("P8").Value = "000010"
$ws.Range("K9").Value = "000000"
$ws.Range("N9").Value = "00000"
$ws.Range("O9").Value = "00000"
$ws.Range("P9").Value = "011000"
$ws.Range("K10").Value = "000000"
$ws.Range("N10").Value = "00000"
$ws.Range("O10").Value = "00000"
$ws.Range("P10").Value = "011001"
$ws.Range("K11").Value = "000000"
$ws.Range("N11").Value = "00000"
$ws.Range("O11").Value = "00000"
$ws.Range("P11").Value = "011010"
$ws.Range("K12").Value = "000000"
$ws.Range("N12").Value = "00000"
$ws.Range("O12").Value = "00000"
$ws.Range("P12").Value = "011011"
$ws.Range("K13").Value = "000000"
$ws.Range("O13").Value = "00000"
$ws.Range("P13").Value = "100100"
$ws.Range("K14").Value = "001100"
$ws.Range("K15").Value = "000000"
$ws.Range("O15").Value = "00000"
$ws.Range("P15").Value = "100101"
$ws.Range("K16").Value = "001101"
$ws.Range("K17").Value = "000000"
$ws.Range("O17").Value = "00000"
$ws.Range("P17").Value = "100110"
$ws.Range("K18").Value = "001110"
$ws.Range("K19").Value = "000000"
$ws.Range("O19").Value = "00000"
$ws.Range("P19").Value = "100111"
$ws.Range("K20").Value = "000000"
$ws.Range("O20").Value = "00000"
$ws.Range("P20").Value = "101010"
$ws.Range("K21").Value = "001010"
$ws.Range("K22").Value = "000000"
$ws.Range("O22").Value = "00000"
$ws.Range("P22").Value = "101011"
$ws.Range("K23").Value = "001011"
$ws.Range("K24").Value = "000000"
$ws.Range("L24").Value = "00000"
$ws.Range("P24").Value = "000000"
$ws.Range("K25").Value = "000000"
$ws.Range("L25").Value = "0000"
$ws.Range("P25").Value = "000010"
$ws.Range("K26").Value = "000000"
$ws.Range("L26").Value = "00000"
$ws.Range("P26").Value = "000011"
$ws.Range("K27").Value = "000000"
$ws.Range("O27").Value = "00000"
$ws.Range("P27").Value = "000100"
$ws.Range("K28").Value = "000000"
$ws.Range("O28").Value = "0000"
$ws.Range("P28").Value = "000110"
$ws.Range("K29").Value = "010011"
$ws.Range("N29").Value = "00000"
$ws.Range("P29").Value = "000101"
$ws.Range("K30").Value = "101011"
$ws.Range("K31").Value = "100000"
$ws.Range("K32").Value = "101000"
$ws.Range("K33").Value = "100001"
$ws.Range("K34").Value = "101001"
$ws.Range("K35").Value = "001111"
$ws.Range("L35").Value = "00000"
$ws.Range("K36").Value = "000000"
$ws.Range("L36").Value = "00000"
$ws.Range("M36").Value = "00000"
$ws.Range("O36").Value = "00000"
$ws.Range("P36").Value = "010000"
$ws.Range("K37").Value = "000000"
$ws.Range("L37").Value = "00000"
$ws.Range("M37").Value = "00000"
$ws.Range("O37").Value = "00000"
$ws.Range("P37").Value = "010010"
$ws.Range("K38").Value = "000010"
$ws.Range("K39").Value = "000100"
$ws.Range("K40").Value = "000101"
$ws.Range("K41").Value = "000001"
$ws.Range("M41").Value = "00001"
$ws.Range("K42").Value = "000111"
$ws.Range("M42").Value = "00000"

# 3) Remove the formatting-only placeholder cells that have no value in the
#    final table (so they do not linger as empty styled cells).
$ws.Range("L2:N2").Clear()
$ws.Range("L3:N3").Clear()
$ws.Range("L4:P4").Clear()
$ws.Range("L5:P5").Clear()
$ws.Range("L6:N6").Clear()
$ws.Range("L7:N7").Clear()
$ws.Range("L8:N8").Clear()
$ws.Range("L9:M9").Clear()
$ws.Range("L10:M10").Clear()
$ws.Range("L11:M11").Clear()
$ws.Range("L12:M12").Clear()
$ws.Range("L13:N13").Clear()
$ws.Range("L14:P14").Clear()
$ws.Range("L15:N15").Clear()
$ws.Range("L16:P16").Clear()
$ws.Range("L17:N17").Clear()
$ws.Range("L18:P18").Clear()
$ws.Range("L19:N19").Clear()
$ws.Range("L20:N20").Clear()
$ws.Range("L21:P21").Clear()
$ws.Range("L22:N22").Clear()
$ws.Range("L23:P23").Clear()
$ws.Range("M24:O24").Clear()
$ws.Range("M25:O25").Clear()
$ws.Range("M26:O26").Clear()
$ws.Range("L27:N27").Clear()
$ws.Range("L28:N28").Clear()
$ws.Range("L29:M29").Clear()
$ws.Range("O29").Clear()
$ws.Range("L30:P30").Clear()
$ws.Range("L31:P31").Clear()
$ws.Range("L32:P32").Clear()
$ws.Range("L33:P33").Clear()
$ws.Range("L34:P34").Clear()
$ws.Range("M35:P35").Clear()
$ws.Range("N36").Clear()
$ws.Range("N37").Clear()
$ws.Range("L38:P38").Clear()
$ws.Range("L39:P39").Clear()
$ws.Range("L40:P40").Clear()
$ws.Range("L41").Clear()
$ws.Range("N41:P41").Clear()
$ws.Range("L42").Clear()
$ws.Range("N42:P42").Clear()

# 4) Leave the selection where the author left it when the file was uploaded.
$ws.Range("N39").Select()
